$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item(1)

# C1/D1 currently share the bold, fully-thin-bordered style (s=1) used by the
# merged header B1:D1. Give them their own plain (non-bold, unaligned) style
# with lighter borders:
#  - C1: top+bottom thin only
#  - D1: top+bottom+right thin only
$c1 = $ws1.Range("C1")
$c1.ClearFormats()
$c1.Borders.Item(8).LineStyle = 1   # xlTop
$c1.Borders.Item(9).LineStyle = 1   # xlBottom

$d1 = $ws1.Range("D1")
$d1.ClearFormats()
$d1.Borders.Item(8).LineStyle = 1   # xlTop
$d1.Borders.Item(10).LineStyle = 1  # xlRight
$d1.Borders.Item(9).LineStyle = 1   # xlBottom

# Rename the "fedcore" column header to "approach"
$ws1.Range("C2").Value2 = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item(2)

# Reuse the exact formats just built on sheet 1 (C1/D1) for sheet 2's
# C1/D1/F1/G1 via copy/paste-format, so the same two new cell styles get
# reused everywhere instead of generating duplicate style-table entries.
$c1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$d1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$c1.Copy()
$ws2.Range("F1").PasteSpecial(-4122)
$d1.Copy()
$ws2.Range("G1").PasteSpecial(-4122)

$ws2.Range("C2").Value2 = "approach"
$ws2.Range("F2").Value2 = "approach"

# Remove the stray empty inline-string cell at G5
$ws2.Range("G5").ClearContents()
